$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata" updates ---
$wsMeta = $wb.Worksheets.Item("Metadata")

# Date: 2026-01-07T21:00:10+00:00 -> 2026-01-14T15:34:52+00:00
$wsMeta.Range("B8").Value = "2026-01-14T15:34:52+00:00"

# Description: "Entrée Resultats d'examens de biologie medicale" -> "Resultats d'examens de biologie medicale"
$wsMeta.Range("B12").Value = "Resultats d'examens de biologie medicale"

# --- Sheet "Elements" updates ---
$wsElem = $wb.Worksheets.Item("Elements")

# Definition for the root element (matches the Description metadata text)
$wsElem.Range("M2").Value = "Resultats d'examens de biologie medicale"

# Short & Definition for fr-lm-resultats-examens-biologie-medicale.laboratoireExecutant
$wsElem.Range("L7").Value = "Laboratoire sous-traitant."
$wsElem.Range("M7").Value = "Laboratoire sous-traitant."

# Short & Definition for fr-lm-resultats-examens-biologie-medicale.auteur
$wsElem.Range("L8").Value = "Participation d'un auteur au document."
$wsElem.Range("M8").Value = "Participation d'un auteur au document."
